# Update the "Skill" value for Sareh Farid (row 7) on Sheet1
# from "System Definition Document" to "Software Requirements Elicitation"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F7").Value = "Software Requirements Elicitation"
